# Insert two new data rows into the "Femacal de La Calera - Zapallo italiano" sheet.
# The new rows are inserted right above the current row 182, pushing every row
# from 182..269 down by two positions (to 184..271), matching the target diff
# (dimension A1:R269 -> A1:R271).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 182 (old row 182 becomes row 184, etc.)
$ws.Rows.Item(182).EntireRow.Insert()
$ws.Rows.Item(182).EntireRow.Insert()

# --- New row 182 ---
$ws.Cells.Item(182, 1).Value = 3
$ws.Cells.Item(182, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(182, 3).Value = "Coquimbo"
$ws.Cells.Item(182, 4).Value = 44523
$ws.Cells.Item(182, 5).Value = 5
$ws.Cells.Item(182, 6).Value = 100112032
$ws.Cells.Item(182, 7).Value = "Zapallo italiano"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 130
$ws.Cells.Item(182, 11).Value = 4000
$ws.Cells.Item(182, 12).Value = 4400
$ws.Cells.Item(182, 13).Value = 4215
$ws.Cells.Item(182, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(182, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(182, 16).Value = 117
$ws.Cells.Item(182, 17).Value = 36
$ws.Cells.Item(182, 18).Value = "Hortaliza"

# --- New row 183 ---
$ws.Cells.Item(183, 1).Value = 3
$ws.Cells.Item(183, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(183, 3).Value = "Coquimbo"
$ws.Cells.Item(183, 4).Value = 44523
$ws.Cells.Item(183, 5).Value = 5
$ws.Cells.Item(183, 6).Value = 100112032
$ws.Cells.Item(183, 7).Value = "Zapallo italiano"
$ws.Cells.Item(183, 8).Value = "Sin especificar"
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 115
$ws.Cells.Item(183, 11).Value = 7000
$ws.Cells.Item(183, 12).Value = 8000
$ws.Cells.Item(183, 13).Value = 7478
$ws.Cells.Item(183, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(183, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(183, 16).Value = 107
$ws.Cells.Item(183, 17).Value = 70
$ws.Cells.Item(183, 18).Value = "Hortaliza"
